$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 294.37036
$ws.Range("I15").Value = 294.37036
$ws.Range("K15").Value = 883.11108
$ws.Range("M15").Value = -714.11108
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H28").Value = 645.6667
$ws.Range("I28").Value = 654.8
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 654.8
$ws.Range("L28").Value = 600
$ws.Range("M28").Value = -169.8
$ws.Range("N28").Value = -1570
$ws.Range("H135").Value = 235.8
$ws.Range("I135").Value = 235.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2122.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 412.7999999999997
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1610.4615
$ws.Range("I137").Value = 1253.6
$ws.Range("K137").Value = 3760.8
$ws.Range("M137").Value = -1210.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4963.8184
$ws.Range("I32").Value = 3189
$ws.Range("K32").Value = 3189
$ws.Range("M32").Value = -2902
$ws.Range("H74").Value = 1403.091
$ws.Range("I74").Value = 1403.091
$ws.Range("K74").Value = 1403.091
$ws.Range("M74").Value = -529.0909999999999
$ws.Range("H77").Value = 1403.091
$ws.Range("I77").Value = 1403.091
$ws.Range("K77").Value = 7015.455
$ws.Range("M77").Value = -2647.455
$ws.Range("H97").Value = 729.3077
$ws.Range("I97").Value = 487.52173
$ws.Range("K97").Value = 487.52173
$ws.Range("M97").Value = 8.478270000000009

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1669.7273
$ws.Range("I31").Value = 1281.5
$ws.Range("J31").Value = 2135.6
$ws.Range("K31").Value = 1281.5
$ws.Range("L31").Value = 2135.6
$ws.Range("M31").Value = -986.5
$ws.Range("N31").Value = -2725.6
$ws.Range("H34").Value = 1669.7273
$ws.Range("I34").Value = 1281.5
$ws.Range("J34").Value = 2135.6
$ws.Range("K34").Value = 1281.5
$ws.Range("L34").Value = 2135.6
$ws.Range("M34").Value = -1079.5
$ws.Range("N34").Value = -2539.6
$ws.Range("H86").Value = 9961488
$ws.Range("I86").Value = 11620736
$ws.Range("K86").Value = 11620736
$ws.Range("M86").Value = -11619613
$ws.Range("H89").Value = 9961488
$ws.Range("I89").Value = 11620736
$ws.Range("K89").Value = 58103680
$ws.Range("M89").Value = -58098064
$ws.Range("H132").Value = 2475.32
$ws.Range("I132").Value = 2705.2273
$ws.Range("K132").Value = 8115.6819
$ws.Range("M132").Value = -5585.6819
$ws.Range("H134").Value = 2262.647
$ws.Range("I134").Value = 1963.8334
$ws.Range("J134").Value = 2979.8
$ws.Range("K134").Value = 5891.5002
$ws.Range("L134").Value = 8939.400000000001
$ws.Range("M134").Value = -3356.5002
$ws.Range("N134").Value = -14009.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 304.7
$ws.Range("I7").Value = 175.25
$ws.Range("J7").Value = 391
$ws.Range("K7").Value = 525.75
$ws.Range("L7").Value = 1173
$ws.Range("M7").Value = -413.75
$ws.Range("N7").Value = -1397
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H34").Value = 547.6667
$ws.Range("I34").Value = 547.6667
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1643.0001
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1559.0001
$ws.Range("N34").ClearContents()
$ws.Range("H55").Value = 3300.75
$ws.Range("I55").Value = 3734.3333
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 11202.9999
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -11025.9999
$ws.Range("N55").Value = -6354
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15880

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 16000
$ws.Range("I20").Value = 16000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 16000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -15755
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 4000000
$ws.Range("I24").Value = 4000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 4000000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3999827
$ws.Range("N24").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 15000
$ws.Range("J13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("N13").Value = -15280
$ws.Range("H16").Value = 415.81818
$ws.Range("I16").Value = 415.81818
$ws.Range("K16").Value = 415.81818
$ws.Range("M16").Value = -245.81818

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 42000
$ws.Range("J80").Value = 42000
$ws.Range("L80").Value = 42000
$ws.Range("N80").Value = -43996
$ws.Range("H83").Value = 42000
$ws.Range("J83").Value = 42000
$ws.Range("L83").Value = 126000
$ws.Range("N83").Value = -135984
$ws.Range("H122").Value = 3509.875
$ws.Range("I122").Value = 3297
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9891
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7441
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 3282.6667
$ws.Range("I126").Value = 3282.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9848.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7378.000100000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2335.682
$ws.Range("I132").Value = 2151.5789
$ws.Range("J132").Value = 3501.6667
$ws.Range("K132").Value = 6454.736699999999
$ws.Range("L132").Value = 10505.0001
$ws.Range("M132").Value = -3924.736699999999
$ws.Range("N132").Value = -15565.0001
